$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 4) with the new projection values
$ws.Range("A4").Value = 86784
$ws.Range("B4").Value = 53958
$ws.Range("C4").Value = 81498
$ws.Range("D4").Value = 4.8213333333333326
$ws.Range("E4").Value = 2.997666666666666
$ws.Range("F4").Value = 4.5276666666666667
$ws.Range("G4").Value = 0.14433333333333334
$ws.Range("H4").Value = 0.1
$ws.Range("I4").Value = 0.12766666666666668
$ws.Range("J4").Value = 0.66811212764870509
$ws.Range("K4").Value = 0.56986671038632719
$ws.Range("L4").Value = 0.60725740832858899
$ws.Range("M4").Value = 4884
$ws.Range("N4").Value = 4608
$ws.Range("O4").Value = 10662
$ws.Range("P4").Value = 0.28395348837209311
$ws.Range("Q4").Value = 0.26790697674418612
$ws.Range("R4").Value = 0.61988372093023258
$ws.Range("S4").Value = 0.0094183333333333341
$ws.Range("T4").Value = 0.0066274999999999997
$ws.Range("U4").Value = 0.024766666666666666
$ws.Range("V4").Value = 0.67135712764870514
$ws.Range("W4").Value = 0.57316171038632713
$ws.Range("X4").Value = 0.61189240832858893

# Remove the now-unused 25th column (Column25) from the table and the sheet
$tbl = $ws.ListObjects.Item(1)
$col = $tbl.ListColumns.Item(25)
$col.Delete()

# Shrink the table back down to the actual data extent (A1:X4)
$tbl.Resize($ws.Range("A1:X4"))
